# Sincronização de dados: atualiza status de um orçamento expirado,
# e adiciona um novo orçamento (PREFEITURA MUNICIPAL DE ROSANA) com seu item.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "quotations": update one status, then insert a new quotation row
# ---------------------------------------------------------------------
$quotations = $wb.Worksheets.Item("quotations")

# Row 2 changed status from Pendente/pending -> Expirada/expired
$quotations.Range("I2").Value = "Expirada"
$quotations.Range("U2").Value = "expired"

# Insert a brand new row at position 15 (pushes existing rows 15-29 down to 16-30)
$quotations.Rows.Item(15).Insert()

# Make sure numeric-looking text values in this new row stay TEXT (matches
# the rest of the sheet, where subtotal/total/discount columns are text).
$quotations.Range("A15:U15").NumberFormat = "@"

$quotations.Range("A15").Value = "NjQxYmI0ZjMtNTE3YS00NjM4LTg0NjktY2Y5ZGExODcxMzc4OjU3MDE2"
$quotations.Range("B15").Value = "SKD1SV5XZX"
$quotations.Range("C15").Value = "PREFEITURA MUNICIPAL DE ROSANA"
$quotations.Range("D15").Value = "Referente a  transformação de 1 maquina composteira no valor total de R`$ 3.000,00, de tensão 380V- trifasico, para 220V/trifasico. `nAs outras 02 são 220v"
$quotations.Range("E15").Value = ""
$quotations.Range("F15").Value = $false
$quotations.Range("G15").Value = "3000"
$quotations.Range("H15").Value = "3000"
$quotations.Range("I15").Value = "Pendente"
$quotations.Range("J15").Value = "2025-11-25T21:47:25.554Z"
$quotations.Range("K15").Value = ""
$quotations.Range("L15").Value = ""
$quotations.Range("M15").Value = "Adriana Vieira Masini"
$quotations.Range("N15").Value = ""
$quotations.Range("O15").Value = "2025-11-18T21:53:35.934Z"
$quotations.Range("P15").Value = "NjYwZjY1NTAtNjFmZS00N2NmLTlmZTktMzY2ZWVjNmViZGJmOjU3MDE2"
$quotations.Range("Q15").Value = "percentage"
$quotations.Range("R15").Value = "0"
$quotations.Range("S15").Value = "0"
$quotations.Range("T15").Value = "NTgwMDkyOTo1NzAxNg=="
$quotations.Range("U15").Value = "pending"

# ---------------------------------------------------------------------
# Sheet "items": insert the matching line item for the new quotation
# ---------------------------------------------------------------------
$items = $wb.Worksheets.Item("items")

# Insert a brand new row at position 63 (pushes existing rows 63-101 down to 64-102)
$items.Rows.Item(63).Insert()

$items.Range("A63").Value = "NjQxYmI0ZjMtNTE3YS00NjM4LTg0NjktY2Y5ZGExODcxMzc4OjU3MDE2"
$items.Range("B63").Value = 1
$items.Range("C63").Value = 3000
$items.Range("D63").Value = "Transformaçao 380V- trifasico, para 220V/trifasico"
$items.Range("E63").Value = 2
$items.Range("F63").Value = "NjQxYmI0ZjMtNTE3YS00NjM4LTg0NjktY2Y5ZGExODcxMzc4OjU3MDE2"
$items.Range("G63").Value = "NDk0NTM0YTgtOGM0MC00ZTYwLWJhNDAtMzMxNDRiZjA3ZGY4OjU3MDE2"
$items.Range("H63").Value = 3000
$items.Range("I63").Value = "service"
$items.Range("J63").Value = "NjQxYmI0ZjMtNTE3YS00NjM4LTg0NjktY2Y5ZGExODcxMzc4OjU3MDE2"
